{"js": "// The whole document is a single paragraph of \"TODO\" notes: three bold\n// section headers (TICKETERA / ATENCION / DISPLAY) followed by free-form\n// lines separated by manual line breaks (<w:br/>), not real paragraph\n// marks. The commit rewrites the notes under each header while leaving\n// the bold headers themselves untouched.\n//\n// Because the new wording does not line up with the old wording run by\n// run (Word's editor commonly splits/merges runs when text is retyped),\n// the most faithful way to reproduce the target markup exactly is to\n// replace the body's single paragraph with the literal OOXML for the\n// post-edit runs, reusing the same bold/size formatting for the three\n// headers that the original document already used.\nconst headerRunProps32 = \"<w:rPr><w:b/><w:bCs/><w:sz w:val=\\\"32\\\"/><w:szCs w:val=\\\"32\\\"/></w:rPr>\";\nconst headerRunProps36 = \"<w:rPr><w:b/><w:bCs/><w:sz w:val=\\\"36\\\"/><w:szCs w:val=\\\"36\\\"/></w:rPr>\";\n\nconst runsXml =\n  \"<w:r>\" + headerRunProps32 + \"<w:t>TICKETERA</w:t></w:r>\" +\n  \"<w:r><w:br/></w:r>\" +\n  \"<w:r><w:t>log ubication, modify query string to ip</w:t></w:r>\" +\n  \"<w:r><w:br/><w:t>c</w:t></w:r>\" +\n  \"<w:r><w:t>heck nums</w:t></w:r>\" +\n  \"<w:r><w:br/></w:r>\" +\n  \"<w:r><w:br/></w:r>\" +\n  \"<w:r>\" + headerRunProps36 + \"<w:t>ATENCION</w:t></w:r>\" +\n  \"<w:r><w:br/></w:r>\" +\n  \"<w:r><w:t>aplicar logs y modificar query</w:t></w:r>\" +\n  \"<w:r><w:br/></w:r>\" +\n  \"<w:r><w:t>relinkear recursos</w:t></w:r>\" +\n  \"<w:r><w:t xml:space=\\\"preserve\\\"> selector</w:t></w:r>\" +\n  \"<w:r><w:br/><w:t xml:space=\\\"preserve\\\"> acomodar ver actuales</w:t></w:r>\" +\n  \"<w:r><w:br/></w:r>\" +\n  \"<w:r><w:t xml:space=\\\"preserve\\\"> ABOUT</w:t></w:r>\" +\n  \"<w:r><w:br/></w:r>\" +\n  \"<w:r><w:br/></w:r>\" +\n  \"<w:r><w:br/></w:r>\" +\n  \"<w:r>\" + headerRunProps36 + \"<w:t>DISPLAY</w:t></w:r>\" +\n  \"<w:r><w:br/></w:r>\" +\n  \"<w:r><w:t>test video logic</w:t></w:r>\" +\n  \"<w:r><w:br/><w:t xml:space=\\\"preserve\\\"> aplicar logs</w:t></w:r>\" +\n  \"<w:r><w:br/><w:t xml:space=\\\"preserve\\\"> modificar query</w:t></w:r>\" +\n  \"<w:r><w:br/><w:t>acomodar barra hora-fecha</w:t></w:r>\" +\n  \"<w:r><w:br/><w:t>cambiar tiempos hora, clima</w:t></w:r>\";\n\nconst ooxml =\n  \"<?xml version=\\\"1.0\\\" encoding=\\\"utf-8\\\" standalone=\\\"yes\\\"?>\" +\n  \"<pkg:package xmlns:pkg=\\\"http://schemas.microsoft.com/office/2006/xmlPackage\\\">\" +\n  \"<pkg:part pkg:name=\\\"/word/document.xml\\\" \" +\n  \"pkg:contentType=\\\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\\\">\" +\n  \"<pkg:xmlData>\" +\n  \"<w:document xmlns:w=\\\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\\\">\" +\n  \"<w:body><w:p>\" + runsXml + \"</w:p></w:body>\" +\n  \"</w:document>\" +\n  \"</pkg:xmlData></pkg:part></pkg:package>\";\n\nconst body = context.document.body;\nbody.insertOoxml(ooxml, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# The whole document is a single paragraph of \"TODO\" notes: three bold\n# section headers (TICKETERA / ATENCION / DISPLAY) followed by free-form\n# lines separated by manual line breaks (<w:br/>), not real paragraph\n# marks. The commit rewrites the notes under each header while leaving\n# the bold headers themselves untouched.\n#\n# Because the new wording does not line up with the old wording run by\n# run (Word's editor commonly splits/merges runs when text is retyped),\n# the most faithful way to reproduce the target markup exactly is to\n# replace the body's single paragraph with the literal OOXML for the\n# post-edit runs, reusing the same bold/size formatting for the three\n# headers that the original document already used.\n\n$headerRunProps32 = '<w:rPr><w:b/><w:bCs/><w:sz w:val=\"32\"/><w:szCs w:val=\"32\"/></w:rPr>'\n$headerRunProps36 = '<w:rPr><w:b/><w:bCs/><w:sz w:val=\"36\"/><w:szCs w:val=\"36\"/></w:rPr>'\n\n$runsXml =\n  ('<w:r>' + $headerRunProps32 + '<w:t>TICKETERA</w:t></w:r>') +\n  '<w:r><w:br/></w:r>' +\n  '<w:r><w:t>log ubication, modify query string to ip</w:t></w:r>' +\n  '<w:r><w:br/><w:t>c</w:t></w:r>' +\n  '<w:r><w:t>heck nums</w:t></w:r>' +\n  '<w:r><w:br/></w:r>' +\n  '<w:r><w:br/></w:r>' +\n  ('<w:r>' + $headerRunProps36 + '<w:t>ATENCION</w:t></w:r>') +\n  '<w:r><w:br/></w:r>' +\n  '<w:r><w:t>aplicar logs y modificar query</w:t></w:r>' +\n  '<w:r><w:br/></w:r>' +\n  '<w:r><w:t>relinkear recursos</w:t></w:r>' +\n  '<w:r><w:t xml:space=\"preserve\"> selector</w:t></w:r>' +\n  '<w:r><w:br/><w:t xml:space=\"preserve\"> acomodar ver actuales</w:t></w:r>' +\n  '<w:r><w:br/></w:r>' +\n  '<w:r><w:t xml:space=\"preserve\"> ABOUT</w:t></w:r>' +\n  '<w:r><w:br/></w:r>' +\n  '<w:r><w:br/></w:r>' +\n  '<w:r><w:br/></w:r>' +\n  ('<w:r>' + $headerRunProps36 + '<w:t>DISPLAY</w:t></w:r>') +\n  '<w:r><w:br/></w:r>' +\n  '<w:r><w:t>test video logic</w:t></w:r>' +\n  '<w:r><w:br/><w:t xml:space=\"preserve\"> aplicar logs</w:t></w:r>' +\n  '<w:r><w:br/><w:t xml:space=\"preserve\"> modificar query</w:t></w:r>' +\n  '<w:r><w:br/><w:t>acomodar barra hora-fecha</w:t></w:r>' +\n  '<w:r><w:br/><w:t>cambiar tiempos hora, clima</w:t></w:r>'\n\n$ooxml =\n  '<?xml version=\"1.0\" encoding=\"utf-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" ' +\n  'pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  ('<w:body><w:p>' + $runsXml + '</w:p></w:body>') +\n  '</w:document>' +\n  '</pkg:xmlData></pkg:part></pkg:package>'\n\n$d = $word.ActiveDocument\n$d.Content.InsertXML($ooxml)\n"}
